$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 152 (A152:E152)
$ws.Cells.Item(152, 1).Value = 151
$ws.Cells.Item(152, 2).Value = 1
$ws.Cells.Item(152, 3).Value = "2024-06-18 04:15:07"
$ws.Cells.Item(152, 4).Value = 200
$ws.Cells.Item(152, 5).Value = 15

# Row 153 (A153:E153)
$ws.Cells.Item(153, 1).Value = 152
$ws.Cells.Item(153, 2).Value = 2
$ws.Cells.Item(153, 3).Value = "2024-06-18 04:15:08"
$ws.Cells.Item(153, 4).Value = 200
$ws.Cells.Item(153, 5).Value = 2
